$d = $word.ActiveDocument

# The empty paragraph that used to follow the "Reliabilty" paragraph is
# removed (merged away), and the "_GoBack" bookmark - which previously sat
# at the very end of the document text (around the last sentence of the
# "Maintenance" paragraph) - is relocated to sit right after the
# "Reliabilty" run instead.

# 1) Remove the empty paragraph right after the "Reliabilty" paragraph by
#    deleting its range (which is just its paragraph mark). This merges it
#    away, leaving the "Reliabilty" paragraph's own paragraph mark (and
#    therefore its paragraph properties) intact.
$emptyPara = $d.Paragraphs.Item(2)
$emptyPara.Range.Delete()

# 2) Move the "_GoBack" bookmark so that it sits right after the
#    "Reliabilty" run (as an empty/collapsed bookmark), instead of at the
#    end of the document. A bookmark placed exactly at the end of a
#    paragraph's text is ambiguous between "end of this paragraph" and
#    "start of the next paragraph", so insert a temporary sentinel
#    character right after the text first; this gives an unambiguous,
#    mid-paragraph insertion point to anchor the bookmark, after which the
#    sentinel is deleted again, leaving the bookmark cleanly between the
#    run and the paragraph mark.
$firstPara = $d.Paragraphs.Item(1)
$textEnd = $firstPara.Range.End - 1

$sentinel = $d.Range($textEnd, $textEnd)
$sentinel.InsertAfter("#")

$bookmarkPoint = $d.Range($textEnd, $textEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)

$sentinelRange = $d.Range($textEnd, $textEnd + 1)
$sentinelRange.Delete()
